$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end and name it "SAYAN"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "SAYAN"

# Fill in header row
$newSheet.Range("A1").Value = "Browser"
$newSheet.Range("B1").Value = "User ID"
$newSheet.Range("C1").Value = "Password"
$newSheet.Range("D1").Value = "Element1"
$newSheet.Range("E1").Value = "Element2"
$newSheet.Range("F1").Value = "Element3"

# Data rows (entered in the same cell order as the original authoring session)
$newSheet.Range("A2").Value = "FF"
$newSheet.Range("A3").Value = "CH"
$newSheet.Range("B3").Value = "HASAN"
$newSheet.Range("B2").Value = "sayan"
$newSheet.Range("C2").Value = "test134"
$newSheet.Range("C3").Value = "yahoo13"
$newSheet.Range("D3").Value = "value1"
$newSheet.Range("D2").Value = "value0"
$newSheet.Range("E2").Value = "value3"
$newSheet.Range("F2").Value = "value5"
$newSheet.Range("E3").Value = "value2"
$newSheet.Range("F3").Value = "vaalue6"

# Apply the same style as the TestCases header row to the new header row
# (bold font + yellow fill, matching existing cellXfs index 1)
$newSheet.Range("A1:F1").Font.Bold = $true
$newSheet.Range("A1:F1").Interior.Color = 65535

# Selections as described in the diff
$wb.Worksheets.Item("TestCases").Range("A2").Select()
$newSheet.Range("L33").Select()

$newSheet.Activate()
